# Atualização automática de preços de eletricidade
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Day (date serial number)
$ws.Range("A2").Value = 45881

# Update hourly prices B2:Z2
$ws.Range("B2").Value = 109.13
$ws.Range("C2").Value = 99.94
$ws.Range("D2").Value = 104.72
$ws.Range("E2").Value = 104.38
$ws.Range("F2").Value = 103
$ws.Range("G2").Value = 104.38
$ws.Range("H2").Value = 99.18000000000001
$ws.Range("I2").Value = 109.7
$ws.Range("J2").Value = 98.69
$ws.Range("K2").Value = 96.13
$ws.Range("L2").Value = 66.34
$ws.Range("M2").Value = 40.94
$ws.Range("N2").Value = 27.2
$ws.Range("O2").Value = 26.28
$ws.Range("P2").Value = 27.2
$ws.Range("Q2").Value = 35.05
$ws.Range("R2").Value = 63.8
$ws.Range("S2").Value = 85.81999999999999
$ws.Range("T2").Value = 96.78
$ws.Range("U2").Value = 121.86
$ws.Range("V2").Value = 154.93
$ws.Range("W2").Value = 140.81
$ws.Range("X2").Value = 126.87
$ws.Range("Y2").Value = 109.7
$ws.Range("Z2").Value = 89.7

# Slot_4h_price
$ws.Range("AB2").Value = 133.08
# Slot_2h_frist_price
$ws.Range("AD2").Value = 147.87
# Slot_2h_second_price
$ws.Range("AF2").Value = 118.28
# Slot_min_price label
$ws.Range("AG2").Value = "10h-17h"
